$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ "B"="12.47380400797446"; "C"="11.49003351170286"; "E"="13.06301130762026"; "F"="16.86991607391245"; "G"="22.68433865537384"; "H"="12.79507591132591"; "I"="19.9257970454268"; "L"="9.944895321098828"; "M"="13.50925961736414"; "N"="17.27861987639285"; "O"="18.6568047516667" }
    3 = @{ "B"="11.99802600893717"; "C"="11.38240386900718"; "E"="13.10988307142015"; "F"="15.89584955866815"; "G"="22.72562814167759"; "H"="12.83880148355165"; "I"="20.0266850126439"; "L"="9.951562977953566"; "M"="13.41335161166051"; "N"="17.31013520986431"; "O"="18.72425339380221" }
    4 = @{ "B"="11.69680833864881"; "C"="11.31555334519014"; "E"="13.14035251059573"; "F"="15.26997757108489"; "G"="22.75964313504788"; "H"="12.86774228358949"; "I"="20.09252285375984"; "L"="9.957018337594677"; "M"="13.35571125051987"; "N"="17.33117255928149"; "O"="18.76998623788411" }
    5 = @{ "B"="11.57194905493191"; "C"="11.28813507493338"; "E"="13.15319484849389"; "F"="15.008197319934"; "G"="22.77567392861717"; "H"="12.8800622004815"; "I"="20.12033130832761"; "L"="9.959584419035821"; "M"="13.33255473622666"; "N"="17.34017032634518"; "O"="18.7897064263958" }
    6 = @{ "B"="11.5510942746371"; "C"="11.28357213309028"; "E"="13.15535305045017"; "F"="14.96433081551589"; "G"="22.77846656609839"; "H"="12.88213970017756"; "I"="20.12500801894072"; "L"="9.96003124736372"; "M"="13.32873022712324"; "N"="17.34169008121707"; "O"="18.79304632902677" }
    7 = @{ "B"="11.69513274353571"; "C"="11.31518426432972"; "E"="13.140523981713"; "F"="15.26647399323133"; "G"="22.75985056101797"; "H"="12.86790630333737"; "I"="20.09289392397929"; "L"="9.957051555159461"; "M"="13.35539758347705"; "N"="17.33129218505649"; "O"="18.77024780666973" }
    8 = @{ "B"="12.31174079459681"; "C"="11.45309145897932"; "E"="13.07882243679675"; "F"="16.5399640634477"; "G"="22.6967722678841"; "H"="12.80971804636787"; "I"="19.95977554324811"; "L"="9.946912179348276"; "M"="13.47594297312209"; "N"="17.28913669544918"; "O"="18.67916320069676" }
    9 = @{ "B"="13.44208439500214"; "C"="11.71674403712448"; "E"="12.97119717137042"; "F"="19.00274580682531"; "G"="22.64214836870593"; "H"="12.71222010044221"; "I"="19.72960143494831"; "L"="9.937800506432897"; "M"="13.72138896475573"; "N"="17.21982264117936"; "O"="18.53492050689511" }
    10 = @{ "B"="14.21675368158522"; "C"="11.90538031655529"; "E"="12.90022320722878"; "F"="20.67494806633232"; "G"="22.64448068034298"; "H"="12.65071493353515"; "I"="19.57929468700943"; "L"="9.937628861146228"; "M"="13.90608762386075"; "N"="17.17699577142696"; "O"="18.45003780088664" }
    11 = @{ "B"="14.55574407043208"; "C"="11.98990146607355"; "E"="12.86968260145963"; "F"="21.3917225636224"; "G"="22.65480444730221"; "H"="12.62493389391979"; "I"="19.5149962337515"; "L"="9.93895619001559"; "M"="13.99080999833205"; "N"="17.15926230537164"; "O"="18.41602991301787" }
    12 = @{ "B"="14.68209218737356"; "C"="12.02170591716328"; "E"="12.85836789575746"; "F"="21.65686569030329"; "G"="22.66004641909764"; "H"="12.61548742873228"; "I"="19.49123434641975"; "L"="9.939659829749102"; "M"="14.02297177251179"; "N"="17.15279785633062"; "O"="18.40381640459896" }
    13 = @{ "B"="14.65497203563313"; "C"="12.0148654937222"; "E"="12.86079359609396"; "F"="21.60004134736742"; "G"="22.65885820889567"; "H"="12.61750782548506"; "I"="19.49632580212326"; "L"="9.939499366636175"; "M"="14.01604198002618"; "N"="17.15417894540696"; "O"="18.40641721205926" }
    14 = @{ "B"="14.56617977441014"; "C"="11.99252217536952"; "E"="12.86874672013202"; "F"="21.4136618050453"; "G"="22.655209003339"; "H"="12.62415038823614"; "I"="19.513029571923"; "L"="9.939010056890757"; "M"="13.99345450292876"; "N"="17.1587254478653"; "O"="18.41501176950575" }
    15 = @{ "B"="14.51152635528995"; "C"="11.97880947999729"; "E"="12.87365082221489"; "F"="21.29868154950795"; "G"="22.65314728913081"; "H"="12.62826034108641"; "I"="19.5233375079799"; "L"="9.938736484128412"; "M"="13.97962869724809"; "N"="17.16154295863656"; "O"="18.4203627915428" }
    16 = @{ "B"="14.1943218707189"; "C"="11.89982929121187"; "E"="12.90225417881469"; "F"="20.62722412089977"; "G"="22.64399248143276"; "H"="12.65244403007964"; "I"="19.58357880273455"; "L"="9.937570306560124"; "M"="13.90056308972107"; "N"="17.17818982941705"; "O"="18.45235316951662" }
    17 = @{ "B"="13.99622302254592"; "C"="11.8510352729396"; "E"="12.92024806260325"; "F"="20.20408069597325"; "G"="22.64074972179023"; "H"="12.6678430275445"; "I"="19.62157925654447"; "L"="9.93721410749267"; "M"="13.8522232800249"; "N"="17.18884957889132"; "O"="18.47315949502004" }
    18 = @{ "B"="13.8810241614674"; "C"="11.82284973585101"; "E"="12.93076203247024"; "F"="19.95656407809801"; "G"="22.63975650537301"; "H"="12.67690699549336"; "I"="19.643819822875"; "L"="9.937141585640839"; "M"="13.82448688158734"; "N"="17.19514541761832"; "O"="18.48556012470007" }
    19 = @{ "B"="13.84180681238173"; "C"="11.81328640271052"; "E"="12.9343501327563"; "F"="19.87204792380568"; "G"="22.63956992708793"; "H"="12.68001142397575"; "I"="19.65141599096149"; "L"="9.937139794713456"; "M"="13.81510804710358"; "N"="17.19730537548047"; "O"="18.48983312799159" }
    20 = @{ "B"="14.01744184960915"; "C"="11.85624207094055"; "E"="12.91831557692848"; "F"="20.24955283636154"; "G"="22.64100466722373"; "H"="12.66618236613522"; "I"="19.61749432879983"; "L"="9.937238334911411"; "M"="13.85736231864973"; "N"="17.18769779605317"; "O"="18.47089975609661" }
    21 = @{ "B"="14.59231570998185"; "C"="11.99909055653253"; "E"="12.86640390555427"; "F"="21.46857628470577"; "G"="22.65624470148373"; "H"="12.62219072195697"; "I"="19.50810734862794"; "L"="9.939148332785482"; "M"="14.00008701469566"; "N"="17.15738322745819"; "O"="18.41246928690029" }
    22 = @{ "B"="14.9562198789134"; "C"="12.09126632653091"; "E"="12.83393560091728"; "F"="22.22866616901552"; "G"="22.6739712907798"; "H"="12.59528310340322"; "I"="19.44003562232081"; "L"="9.941567696680851"; "M"="14.09381772157382"; "N"="17.13903263970085"; "O"="18.37815600280654" }
    23 = @{ "B"="14.76310502877205"; "C"="12.04218400820078"; "E"="12.85113127345216"; "F"="21.82633154458858"; "G"="22.66379991613414"; "H"="12.60947546611706"; "I"="19.47605383396446"; "L"="9.940169665040834"; "M"="14.04375770305492"; "N"="17.14869315294709"; "O"="18.3961144294617" }
    24 = @{ "B"="14.00785288863796"; "C"="11.85388849124389"; "E"="12.91918872709464"; "F"="20.22900810905287"; "G"="22.64088669310504"; "H"="12.66693249429789"; "I"="19.61933989862282"; "L"="9.937226969636653"; "M"="13.85503878692751"; "N"="17.18821799556755"; "O"="18.47192001742074" }
    25 = @{ "B"="13.14559944123111"; "C"="11.64624246136419"; "E"="12.99888671096407"; "F"="18.34778573295695"; "G"="22.64948776594514"; "H"="12.7368176454013"; "I"="19.78856670991416"; "L"="9.939116496914533"; "M"="13.65413739703139"; "N"="17.23714877715621"; "O"="18.57024746162893" }
}

foreach ($row in $data.Keys) {
    $rowVals = $data[$row]
    foreach ($col in $rowVals.Keys) {
        $ws.Range("$col$row").Value = [double]$rowVals[$col]
    }
}
